# Target worksheet: "8.1" (xl/worksheets/sheet2.xml)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("8.1")

# Update WORKORDER (column E) values for the rows that stay in the sheet
$ws.Range("E3").Value = 191020929
$ws.Range("E4").Value = 191020946
$ws.Range("E5").Value = 191020871
$ws.Range("E6").Value = 191020870
$ws.Range("E7").Value = 191020925
$ws.Range("E9").Value = 191023017
$ws.Range("E10").Value = 191022896
$ws.Range("E11").Value = 191022897

# Remove rows 12 through 22 entirely (shrinks dimension to A1:K11)
$ws.Rows("12:22").Delete()
